# Applies the "Testing for new Fibonacci program" edit:
#   1) Refresh the cached "datetimeFigureOut" date stamp (2024/2/15 -> 2024/2/28)
#      on the slide master and every slide layout's Date placeholder.
#   2) On slide 2's opcode-table text box, split the carry bit between the
#      "010" and "110" rows:
#        010:Rd=R0+R1+Carry  ->  010:Rd=R0+R1
#        110:Rd=R0+R1        ->  110:Rd=R0  +R1+Carry

$p = $ppt.ActivePresentation

# --- 1) Date placeholder (ppPlaceholderDate = 16) cached text ---
function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) { $isDatePh = $true }
        } catch {
        }
        if ($isDatePh -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $len = $tr.Length
            if ($len -gt 0) {
                $full = $tr.Characters(1, $len)
                $full.Text = $newText
            }
        }
    }
}

Set-DatePlaceholderText $p.SlideMaster.Shapes "2024/2/28"
for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($L)
    Set-DatePlaceholderText $layout.Shapes "2024/2/28"
}

# --- 2) Slide 2 opcode table text box ("文本框 18", shape 13) ---
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(13)
$tr = $shp.TextFrame.TextRange

# paragraph 4: "010:Rd=R0+R1+Carry" -> "010:Rd=R0+R1"
$para4 = $tr.Paragraphs(4, 1)
$para4.Text = "010:Rd=R0+R1"

# paragraph 8: "110:Rd=R0+R1" -> "110:Rd=R0" followed by a separate "+R1+Carry" run
$para8 = $tr.Paragraphs(8, 1)
$tail = $para8.Characters(10, 3)
$tail.Text = "+R1+Carry"
